$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "A queen sat at a window."
$ws.Range("C2").Value = "queen, window, sit"

$ws.Range("B3").Value = "A queen had a wonderful looking-glass."
$ws.Range("C3").Value = "queen, looking-glass"

$ws.Range("B4").Value = "When the queen stood in front of looking-glass and looked at herself in it, and said, " + [char]8220 + "looking-glass, looking-glass, who in this land is the fairest of all." + [char]8221
$ws.Range("C4").Value = "queen, looking-glass, say, fairest"

$ws.Range("B5").Value = "The looking-glass answered, " + [char]8220 + "Queen, you the fairest of all." + [char]8221
$ws.Range("C5").Value = "looking-glass, answer, fairest"

$ws.Range("B6").Value = "Snow-white was talking with animal friends."
$ws.Range("C6").Value = "Snow-white, talk, animal"

$ws.Range("B7").Value = "But now the poor Snow-white was all alone in the great forest."
$ws.Range("C7").Value = "Snow-white, forest"

$ws.Range("B8").Value = "Show-white ran as long as her feet would go until it was almost evening, then she saw a little cottage and went into it to rest herself. "
$ws.Range("C8").Value = "Snow-white, run, cottage, rest"

$ws.Range("B9").Value = "Little snow-white was so hungry and thirsty that she ate some vegetables and bread from plate and drank water out of mug."
$ws.Range("C9").Value = "Snow-white, hungry, thirsty, ate vegetables, bread, plate, mug, water"

$ws.Range("B10").Value = "As Snow-white was so tired, she laid herself down on one of the little beds, and went to sleep."
$ws.Range("C10").Value = "Snow-white, tired, lay, bed, sleep"

$ws.Range("B11").Value = "When it was morning little Snow-white awoke, and was frightened when she saw the seven dwarfs."
$ws.Range("C11").Value = "mornig, Snow-white, awoke, frightened, seven-dwarfs"

$ws.Range("B12").Value = "But the seven dwarfs were friendly and asked her what her name was.  My name is snow-white, she answered.  "
$ws.Range("C12").Value = "seven-dwarfs, snow-white, name, ask, answer"

$ws.Range("B13").Value = "The queen went into a quite secret, lonely room, where no one ever came, and there she made a very poisonous apple."
$ws.Range("C13").Value = "queen, secret, room, made, apple"

$win = $excel.ActiveWindow
$ws.Range("B17").Select()
$win.ScrollRow = 1
$win.ScrollColumn = 2
